$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FRCcombined column (W) values for rows 2-22. These numbers
# come from re-running the FRC analysis on renamed/additional raw data
# files (see commit message) -- only the computed values change, the
# sheet layout/headers stay the same.
$ws.Range("W2").Value  = 10.362694300518134
$ws.Range("W3").Value  = 11.049723756906078
$ws.Range("W4").Value  = 10.526315789473683
$ws.Range("W5").Value  = 44.444444444444443
$ws.Range("W6").Value  = 10.695187165775401
$ws.Range("W7").Value  = 6.8259385665529004
$ws.Range("W8").Value  = 7.9365079365079358
$ws.Range("W9").Value  = 7.4074074074074066
$ws.Range("W10").Value = 6.4935064935064926
$ws.Range("W11").Value = 6.8728522336769755
$ws.Range("W12").Value = 7.7220077220077208
$ws.Range("W13").Value = 8.3333333333333321
$ws.Range("W14").Value = 6.4516129032258061
$ws.Range("W15").Value = 6.9204152249134943
$ws.Range("W16").Value = 6.3897763578274756
$ws.Range("W17").Value = 7.3529411764705879
$ws.Range("W18").Value = 6.9930069930069925
$ws.Range("W19").Value = 5.9523809523809526
$ws.Range("W20").Value = 6.6225165562913899
$ws.Range("W21").Value = 40
$ws.Range("W22").Value = 166.66666666666666
